# Generate Report for Handback
# Update the timestamp text values recorded on the report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 19:09:35"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 19:09:30"
$wsZhCn.Range("K2").Value = "2016-08-23 19:09:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
# de-de!H2 shares its original shared-string value with Overview!G2
# ("2016-08-23 19:08:41"), so it moves in lockstep with that update.
$wsDeDe.Range("H2").Value = "2016-08-23 19:09:35"
$wsDeDe.Range("K2").Value = "2016-08-23 19:09:55"
